# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect freshly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Row => [oldValue, newValue] mapping for the "展览" sheet (F column)
$sheetExhibitionUpdates = @{
    3  = 5259
    4  = 13
    5  = 7566
    6  = 47
    12 = 4366
    13 = 1778
    14 = 111
    16 = 2949
    18 = 569
    19 = 215
    20 = 528
    21 = 458
    23 = 326
    25 = 1709
    26 = 1215
    27 = 96
    28 = 1407
    29 = 115
    30 = 588
    31 = 32
    32 = 517
    34 = 16
    37 = 72
    38 = 3005
    39 = 712
    40 = 36
    41 = 113
    43 = 65
}

# Row => newValue mapping for the "全部类型" sheet (F column)
$sheetAllTypesUpdates = @{
    3  = 5259
    4  = 13
    5  = 7566
    6  = 47
    12 = 4366
    13 = 1778
    14 = 111
    16 = 2949
    18 = 569
    19 = 215
    20 = 528
    21 = 458
    24 = 326
    26 = 1709
    27 = 1215
    28 = 96
    29 = 1407
    30 = 115
    31 = 588
    32 = 32
    35 = 16
    38 = 72
    39 = 3005
    41 = 712
    42 = 36
    43 = 113
    45 = 65
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $sheetExhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $sheetExhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheetAllTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $sheetAllTypesUpdates[$row]
}
